$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7572.8213
$ws.Range("J19").Value = 12237.529
$ws.Range("L19").Value = 12237.529
$ws.Range("N19").Value = -12587.529
$ws.Range("H40").Value = 3436.5
$ws.Range("I40").Value = 1433.3334
$ws.Range("J40").Value = 4638.4
$ws.Range("K40").Value = 1433.3334
$ws.Range("L40").Value = 4638.4
$ws.Range("M40").Value = -1258.3334
$ws.Range("N40").Value = -4988.4
$ws.Range("H43").Value = 1510.871
$ws.Range("I43").Value = 1272.1
$ws.Range("J43").Value = 1624.5714
$ws.Range("K43").Value = 1272.1
$ws.Range("L43").Value = 1624.5714
$ws.Range("M43").Value = -1203.1
$ws.Range("N43").Value = -1762.5714
$ws.Range("H64").Value = 2950
$ws.Range("I64").Value = 2866.6667
$ws.Range("K64").Value = 2866.6667
$ws.Range("M64").Value = -2618.6667
$ws.Range("H67").Value = 2950
$ws.Range("I67").Value = 2866.6667
$ws.Range("K67").Value = 2866.6667
$ws.Range("M67").Value = -2008.6667
$ws.Range("H74").Value = 3176.3157
$ws.Range("I74").Value = 3239.2856
$ws.Range("K74").Value = 3239.2856
$ws.Range("M74").Value = -2303.2856
$ws.Range("H76").Value = 3312.5
$ws.Range("I76").Value = 3350
$ws.Range("K76").Value = 3350
$ws.Range("M76").Value = -3035
$ws.Range("H77").Value = 3176.3157
$ws.Range("I77").Value = 3239.2856
$ws.Range("K77").Value = 16196.428
$ws.Range("M77").Value = -11516.428
$ws.Range("H79").Value = 3312.5
$ws.Range("I79").Value = 3350
$ws.Range("K79").Value = 3350
$ws.Range("M79").Value = -2258
$ws.Range("H132").Value = 3869.3572
$ws.Range("I132").Value = 3667.0386
$ws.Range("K132").Value = 11001.1158
$ws.Range("M132").Value = -8471.1158
$ws.Range("H138").Value = 3215.04
$ws.Range("I138").Value = 1479.7142
$ws.Range("J138").Value = 3497.535
$ws.Range("K138").Value = 4439.142599999999
$ws.Range("L138").Value = 10492.605
$ws.Range("M138").Value = 700.8574000000008
$ws.Range("N138").Value = -20772.605
$ws.Range("H140").Value = 74178.56
$ws.Range("J140").Value = 74178.56
$ws.Range("L140").Value = 74178.56
$ws.Range("N140").Value = -84538.56

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 266.66666
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H32").Value = 6858096.5
$ws.Range("I32").Value = 7150943.5
$ws.Range("K32").Value = 7150943.5
$ws.Range("M32").Value = -7150656.5
$ws.Range("H61").Value = 11497745
$ws.Range("I61").Value = 19609872
$ws.Range("J61").Value = 5566.6665
$ws.Range("K61").Value = 19609872
$ws.Range("L61").Value = 5566.6665
$ws.Range("M61").Value = -19609660
$ws.Range("N61").Value = -5990.6665
$ws.Range("H63").Value = 3599.0386
$ws.Range("I63").Value = 2592.0715
$ws.Range("J63").Value = 4773.8335
$ws.Range("K63").Value = 2592.0715
$ws.Range("L63").Value = 4773.8335
$ws.Range("M63").Value = -1906.0715
$ws.Range("N63").Value = -6145.8335
$ws.Range("H66").Value = 3599.0386
$ws.Range("I66").Value = 2592.0715
$ws.Range("J66").Value = 4773.8335
$ws.Range("K66").Value = 12960.3575
$ws.Range("L66").Value = 23869.1675
$ws.Range("M66").Value = -9528.3575
$ws.Range("N66").Value = -30733.1675
$ws.Range("H88").Value = 2975
$ws.Range("I88").Value = 2633.3333
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 2633.3333
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -2227.3333
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 2975
$ws.Range("I91").Value = 2633.3333
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 2633.3333
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -1229.3333
$ws.Range("N91").Value = -6808
$ws.Range("H122").Value = 60391.53
$ws.Range("I122").Value = 84670.086
$ws.Range("J122").Value = 2123
$ws.Range("K122").Value = 254010.258
$ws.Range("L122").Value = 6369
$ws.Range("M122").Value = -251560.258
$ws.Range("N122").Value = -11269
$ws.Range("H132").Value = 1511492
$ws.Range("I132").Value = 2449.5278
$ws.Range("K132").Value = 7348.5834
$ws.Range("M132").Value = -4818.5834
$ws.Range("H136").Value = 11497745
$ws.Range("I136").Value = 19609872
$ws.Range("J136").Value = 5566.6665
$ws.Range("K136").Value = 58829616
$ws.Range("L136").Value = 16699.9995
$ws.Range("M136").Value = -58827066
$ws.Range("N136").Value = -21799.9995

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 266.66666
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H22").Value = 576.5714
$ws.Range("I22").Value = 760
$ws.Range("J22").Value = 118
$ws.Range("K22").Value = 760
$ws.Range("L22").Value = 118
$ws.Range("M22").Value = -587
$ws.Range("N22").Value = -464
$ws.Range("H39").Value = 79026.5
$ws.Range("J39").Value = 79026.5
$ws.Range("L39").Value = 79026.5
$ws.Range("N39").Value = -79804.5
$ws.Range("H86").Value = 1819.8
$ws.Range("I86").Value = 1905.1578
$ws.Range("J86").Value = 1549.5
$ws.Range("K86").Value = 1905.1578
$ws.Range("L86").Value = 1549.5
$ws.Range("M86").Value = -782.1578
$ws.Range("N86").Value = -3795.5
$ws.Range("H89").Value = 1819.8
$ws.Range("I89").Value = 1905.1578
$ws.Range("J89").Value = 1549.5
$ws.Range("K89").Value = 9525.789000000001
$ws.Range("L89").Value = 7747.5
$ws.Range("M89").Value = -3909.789000000001
$ws.Range("N89").Value = -18979.5
$ws.Range("H134").Value = 2555.6099
$ws.Range("I134").Value = 2365.7144
$ws.Range("J134").Value = 3663.3333
$ws.Range("K134").Value = 7097.1432
$ws.Range("L134").Value = 10989.9999
$ws.Range("M134").Value = -4562.1432
$ws.Range("N134").Value = -16059.9999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4011.1562
$ws.Range("I62").Value = 3866.5
$ws.Range("J62").Value = 4197.143
$ws.Range("K62").Value = 3866.5
$ws.Range("L62").Value = 4197.143
$ws.Range("M62").Value = -3242.5
$ws.Range("N62").Value = -5445.143
$ws.Range("H65").Value = 4011.1562
$ws.Range("I65").Value = 3866.5
$ws.Range("J65").Value = 4197.143
$ws.Range("K65").Value = 19332.5
$ws.Range("L65").Value = 20985.715
$ws.Range("M65").Value = -16212.5
$ws.Range("N65").Value = -27225.715
$ws.Range("H132").Value = 30305600
$ws.Range("I132").Value = 38463772
$ws.Range("K132").Value = 115391316
$ws.Range("M132").Value = -115388786
$ws.Range("H135").Value = 50097.145
$ws.Range("J135").Value = 50097.145
$ws.Range("L135").Value = 50097.145
$ws.Range("N135").Value = -60237.145
$ws.Range("H137").Value = 49999.5
$ws.Range("J137").Value = 49999.5
$ws.Range("L137").Value = 49999.5
$ws.Range("N137").Value = -60199.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1989.55
$ws.Range("I107").Value = 374.875
$ws.Range("J107").Value = 3066
$ws.Range("K107").Value = 1124.625
$ws.Range("L107").Value = 9198
$ws.Range("M107").Value = 795.375
$ws.Range("N107").Value = -13038

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 23333.334
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = -9064
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 23333.334
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 150000
$ws.Range("M77").Value = -25320
$ws.Range("N77").Value = -159360
$ws.Range("H80").Value = 1881740
$ws.Range("I80").Value = 3001900
$ws.Range("J80").Value = 201500
$ws.Range("K80").Value = 3001900
$ws.Range("L80").Value = 201500
$ws.Range("M80").Value = -3000902
$ws.Range("N80").Value = -203496
$ws.Range("H83").Value = 1881740
$ws.Range("I83").Value = 3001900
$ws.Range("J83").Value = 201500
$ws.Range("K83").Value = 15009500
$ws.Range("L83").Value = 1007500
$ws.Range("M83").Value = -15004508
$ws.Range("N83").Value = -1017484
$ws.Range("H122").Value = 2173.2727
$ws.Range("I122").Value = 2111.7778
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 6335.3334
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -3885.3334
$ws.Range("N122").Value = -12250
$ws.Range("H134").Value = 18141.834
$ws.Range("J134").Value = 18141.834
$ws.Range("L134").Value = 54425.50199999999
$ws.Range("N134").Value = -59495.50199999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5001
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 10480.3
$ws.Range("I22").Value = 674
$ws.Range("J22").Value = 20286.6
$ws.Range("K22").Value = 674
$ws.Range("L22").Value = 20286.6
$ws.Range("M22").Value = -379
$ws.Range("N22").Value = -20876.6
$ws.Range("H27").Value = 10480.3
$ws.Range("I27").Value = 674
$ws.Range("J27").Value = 20286.6
$ws.Range("K27").Value = 674
$ws.Range("L27").Value = 20286.6
$ws.Range("M27").Value = -567
$ws.Range("N27").Value = -20500.6
$ws.Range("H122").Value = 5537.1562
$ws.Range("I122").Value = 3775
$ws.Range("J122").Value = 5943.8076
$ws.Range("K122").Value = 11325
$ws.Range("L122").Value = 17831.4228
$ws.Range("M122").Value = -8875
$ws.Range("N122").Value = -22731.4228
